$d = $word.ActiveDocument

# The document contains six "<id>p029v_N</id>" markers (N = 1..6), each
# split across three runs: one run for "<id>" (Courier New / 7f6000 / 18pt),
# one run for the bare "p029v_N" text (Arial / 000000 / 22pt), and one run
# for "</id>" (Courier New / 7f6000 / 18pt). The edit collapses each triple
# into a single run carrying the formatting of the first ("<id>") run,
# leaving the rendered text unchanged. A seventh, visually identical
# "<id>fig_p029v_1</id>" marker exists elsewhere in the document and must
# stay untouched (it is not part of this change).

for ($n = 1; $n -le 6; $n++) {
    $needle = "<id>p029v_$n</id>"
    $rng = $d.Content
    $rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
